# Update RW_efficacy_in_%, RW_lower, RW_upper (columns Z, AA, AB) for the
# rows whose random-effects meta-analysis recalculation changed the
# resulting values (PythonMeta test code).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 12; Z = 94.43796319779074;  AA = -3.302327052537257;  AB = 99.70052704259608 },
    @{ Row = 23; Z = 97.15223914397328;  AA = 50.66493499150665;   AB = 99.83561911002404 },
    @{ Row = 26; Z = 83.86268848251524;  AA = -296.1118075190292;  AB = 99.34257747921377 },
    @{ Row = 27; Z = 96.97505086197665;  AA = 49.63835990348066;   AB = 99.8183077979571  },
    @{ Row = 30; Z = 85.72256728778468;  AA = -176.3903643897301;  AB = 99.2624739820365  },
    @{ Row = 36; Z = 94.5169636319258;   AA = 5.021905453613229;   AB = 99.68346714095279 },
    @{ Row = 46; Z = 99.2008996601827;   AA = 86.7896061036983;    AB = 99.95166220189128 },
    @{ Row = 49; Z = 95.31787684795356;  AA = 20.11694832456909;   AB = 99.72557035877897 },
    @{ Row = 50; Z = 92.43657029284805;  AA = -34.22888656146353;  AB = 99.57382147464335 },
    @{ Row = 56; Z = 92.07499501362942;  AA = -40.6707088182735;   AB = 99.5535267820742  },
    @{ Row = 57; Z = 93.93970206924602;  AA = -4.991807763674427;  AB = 99.65018974535451 }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("Z$r").Value = $u.Z
    $ws.Range("AA$r").Value = $u.AA
    $ws.Range("AB$r").Value = $u.AB
}
